$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New time-log entry for row 84 (previously blank).
$ws.Range("A84").Value = 41955
$ws.Range("B84").Value = 0.87708333333333333
$ws.Range("C84").Value = 0.96250000000000002
$ws.Range("D84").Value = 30
$ws.Range("F84").Value = "Coding"

# Re-assert the (unchanged) shared formula on E84 so the engine
# re-evaluates it against the new B84/C84/D84 inputs instead of
# keeping a stale cached result for this previously-blank row.
$ws.Range("E84").Formula = '=IF(AND(NOT(ISBLANK(B84)),NOT(ISBLANK(C84))),(C84-B84)*24-D84/60,"")'

$excel.CalculateFull()

# Matches the author's final cursor position recorded in the sheet view.
$ws.Range("C85").Select()
